# Update "Datos actualizados" timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 20:38"

# Swap Camerun/Libano rows 79 and 80: Libano moves to row 79 (with updated stats),
# Camerun moves to row 80 (with its previous, unchanged stats).
$ws.Range("A79").Value = "Libano"
$ws.Range("B79").Value = 20011
$ws.Range("C79").Value = 521
$ws.Range("D79").Value = 5868
$ws.Range("E79").Value = 13956
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 187

$ws.Range("A80").Value = "Camerun"
$ws.Range("B80").Value = 19604
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 18448
$ws.Range("E80").Value = 741
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 415

# Row 4 (Estados Unidos) updates
$ws.Range("B4").Value = 6410295
$ws.Range("C4").Value = 21238
$ws.Range("D4").Value = 3643636
$ws.Range("E4").Value = 2574234
$ws.Range("G4").Value = 314
$ws.Range("H4").Value = 192425

# Row 5 (India) updates
$ws.Range("B5").Value = 4109476
$ws.Range("C5").Value = 89237
$ws.Range("D5").Value = 3177667
$ws.Range("E5").Value = 861136
$ws.Range("G5").Value = 1038
$ws.Range("H5").Value = 70673

# Row 24 (Alemania) updates
$ws.Range("B24").Value = 250791
$ws.Range("C24").Value = 510
$ws.Range("E24").Value = 15182

# Row 32 (Ecuador) updates
$ws.Range("B32").Value = 118045
$ws.Range("C32").Value = 870
$ws.Range("D32").Value = 102304
$ws.Range("E32").Value = 9017
$ws.Range("G32").Value = 50
$ws.Range("H32").Value = 6724

# Row 97 (Guayana Francesa) updates
$ws.Range("B97").Value = 9322
$ws.Range("C97").Value = 46
$ws.Range("D97").Value = 8870
$ws.Range("E97").Value = 390

# Row 100 (Maldivas) updates
$ws.Range("B100").Value = 8486
$ws.Range("C100").Value = 125
$ws.Range("D100").Value = 5822
$ws.Range("E100").Value = 2635

$wb.Save()
